# Edit script: insert 3 new data rows (weekly "Sandia" price records) at the
# top of the data block (row 35) on the single worksheet, shifting all
# existing records down by 3 rows. Dimension grows from A1:R123 to A1:R126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 35; everything currently at
# row 35 and below (through 123) moves down to rows 38-126.
$ws.Rows("35:37").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$categoriaId = 100112028
$categoria   = "Sandia"
$variedad    = "Sin especificar"
$kgOUnidades = 1
$clasif      = "Hortaliza"

# New rows' specific data (Fecha, Calidad, Volumen, PrecioMin, PrecioMax,
# PrecioPromedioPonderado, UnidadComercializacion, Origen, Precio$/Kg).
$newRows = @(
    @{ Row=35; Fecha=44560; Calidad="Extra";   J=500; K=2500; L=2500; M=2500; N="`$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=2500 },
    @{ Row=36; Fecha=44560; Calidad="Primera"; J=600; K=2000; L=2200; M=2100; N="`$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=2100 },
    @{ Row=37; Fecha=44560; Calidad="Segunda"; J=200; K=1800; L=1800; M=1800; N="`$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=1800 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $categoriaId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $kgOUnidades
    $ws.Cells.Item($row, 18).Value = $clasif
}
